# fix: selenium pipeline in wireframe generator
#
# The "Effort Estimation" sheet's feature breakdown changes:
#   - "Login" -> "User Login"
#   - "Filtering and Sorting" -> split into "Product Search" + new "Pagination" row
#   - "Add to Cart" -> "Reviews and Ratings"
#   - a new "Password Recovery" row is added under "User Registration"
# and the per-row effort numbers / cost-summary totals are recomputed to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Effort Estimation")
$ws2 = $wb.Worksheets.Item("Cost Summary")

# --- Sheet1: make room for the two brand-new subfeature rows -----------------
# (insert the lower one first so the earlier insert point doesn't shift)
$ws1.Rows("8:8").Insert()
$ws1.Rows("4:4").Insert()

# --- Sheet1: rewrite every data row (2-11) to match the new breakdown --------
$rows = @(
    @("User Authentication", "User Registration", "Password Recovery",      "5", "1",                  "0.8999999999999999", "3", "0.6000000000000001", "0.54"),
    @("User Authentication", "User Registration", "Email Verification",     "3", "0.6000000000000001", "0.54",                "2", "0.4",                "0.36"),
    @("User Authentication", "User Registration", "Password Recovery",      "4", "0.8",                 "0.72",               "3", "0.6000000000000001", "0.54"),
    @("User Authentication", "User Login",         "Frontend Implementation","3","0.6000000000000001",  "0.54",               "2", "0.4",                "0.36"),
    @("User Authentication", "User Login",         "Session Management",    "2", "0.4",                 "0.36",               "2", "0.4",                "0.36"),
    @("Product Catalog",     "Product Listing",     "Frontend Implementation","7","1.4",                 "1.26",               "5", "1",                  "0.8999999999999999"),
    @("Product Catalog",     "Product Listing",     "Product Search",       "5", "1",                   "0.8999999999999999", "4", "0.8",                "0.72"),
    @("Product Catalog",     "Product Listing",     "Pagination",           "2", "0.4",                 "0.36",               "2", "0.4",                "0.36"),
    @("Product Catalog",     "Product Details",      "Frontend Implementation","5","1",                  "0.8999999999999999", "3", "0.6000000000000001", "0.54"),
    @("Product Catalog",     "Product Details",      "Reviews and Ratings",  "4", "0.8",                 "0.72",               "3", "0.6000000000000001", "0.54")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r   = $i + 2
    $row = $rows[$i]
    $ws1.Cells.Item($r, 1).Value2 = $row[0]
    $ws1.Cells.Item($r, 2).Value2 = $row[1]
    $ws1.Cells.Item($r, 3).Value2 = $row[2]
    $ws1.Cells.Item($r, 4).Value2 = [double]$row[3]
    $ws1.Cells.Item($r, 5).Value2 = [double]$row[4]
    $ws1.Cells.Item($r, 6).Value2 = [double]$row[5]
    $ws1.Cells.Item($r, 7).Value2 = [double]$row[6]
    $ws1.Cells.Item($r, 8).Value2 = [double]$row[7]
    $ws1.Cells.Item($r, 9).Value2 = [double]$row[8]
}

# --- Sheet1: Total row (12) ---------------------------------------------------
$ws1.Cells.Item(12, 1).Value2 = "Total"
$ws1.Cells.Item(12, 3).Value2 = "Total"
$ws1.Cells.Item(12, 4).Value2 = [double]40
$ws1.Cells.Item(12, 5).Value2 = [double]8
$ws1.Cells.Item(12, 6).Value2 = [double]7.2
$ws1.Cells.Item(12, 7).Value2 = [double]29
$ws1.Cells.Item(12, 8).Value2 = [double]"5.799999999999999"
$ws1.Cells.Item(12, 9).Value2 = [double]5.22

# --- Sheet1: Units row (13) ---------------------------------------------------
$ws1.Cells.Item(13, 3).Value2 = "Units"
$ws1.Cells.Item(13, 4).Value2 = "days"
$ws1.Cells.Item(13, 5).Value2 = "days"
$ws1.Cells.Item(13, 6).Value2 = "days"
$ws1.Cells.Item(13, 7).Value2 = "days"
$ws1.Cells.Item(13, 8).Value2 = "days"
$ws1.Cells.Item(13, 9).Value2 = "days"

# --- Sheet2 (Cost Summary): updated pricing ----------------------------------
$ws2.Cells.Item(2, 2).Value2 = [double]52.8
$ws2.Cells.Item(2, 4).Value2 = "₹6,336.00"

$ws2.Cells.Item(3, 2).Value2 = [double]38.28
$ws2.Cells.Item(3, 4).Value2 = "₹4,899.84"

$ws2.Cells.Item(4, 2).Value2 = [double]"7.920000000000001"
$ws2.Cells.Item(4, 4).Value2 = "₹760.32"

$ws2.Cells.Item(5, 4).Value2 = "₹11,996.16"

# --- Sheet2: column B widened by one character (19.71 -> 20.71 chars) --------
$ws2.Columns("B:B").ColumnWidth = 20
